$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.09541473491200993
$ws.Cells.Item(2, 3).Value = 0.00000005542031828998138
$ws.Cells.Item(2, 4).Value = 0.000000006899374246642727

$ws.Cells.Item(3, 2).Value = 0.0629129249586846
$ws.Cells.Item(3, 3).Value = 0.000000206879638728754
$ws.Cells.Item(3, 4).Value = 0.00000002684068627588077

$ws.Cells.Item(4, 2).Value = 0.1355369631985098
$ws.Cells.Item(4, 3).Value = 0.0000005420385614819802
$ws.Cells.Item(4, 4).Value = 0.00000004718708395286664

$ws.Cells.Item(5, 2).Value = 0.1426967817146367
$ws.Cells.Item(5, 3).Value = 0.000001166672082080816
$ws.Cells.Item(5, 4).Value = 0.0000001432830778975021

$ws.Cells.Item(6, 2).Value = 0.1480364181349707
$ws.Cells.Item(6, 3).Value = 0.000002084659349794333
$ws.Cells.Item(6, 4).Value = 0.0000001927824229906775

$ws.Cells.Item(7, 2).Value = 0.1490360290672124
$ws.Cells.Item(7, 3).Value = 0.000003866901263860968
$ws.Cells.Item(7, 4).Value = 0.0000004519348725684993

$ws.Cells.Item(8, 2).Value = 0.1489651913653108
$ws.Cells.Item(8, 3).Value = 0.000006173437684280043
$ws.Cells.Item(8, 4).Value = 0.0000006559406686853016

$ws.Cells.Item(9, 2).Value = 0.1002332623828828
$ws.Cells.Item(9, 3).Value = 0.000008608163214611935
$ws.Cells.Item(9, 4).Value = 0.0000006138005831898286

$ws.Cells.Item(10, 2).Value = 0.08542122415966777
$ws.Cells.Item(10, 3).Value = 0.0000135788108587008
$ws.Cells.Item(10, 4).Value = 0.000001583333715730562

$ws.Cells.Item(11, 2).Value = 0.1057106568782109
$ws.Cells.Item(11, 3).Value = 0.00002069694003600734
$ws.Cells.Item(11, 4).Value = 0.000001530111086025295

$ws.Cells.Item(12, 2).Value = 0.122438029519025
$ws.Cells.Item(12, 3).Value = 0.00002514263863981077
$ws.Cells.Item(12, 4).Value = 0.000001733408530084016

$ws.Cells.Item(13, 2).Value = 0.1056591056193066
$ws.Cells.Item(13, 3).Value = 0.00003288472981545815
$ws.Cells.Item(13, 4).Value = 0.000003262770660277787

